$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldVal = "System, dnasr281@gmail.com"
$newVal = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $oldVal) {
        $cell.Value = $newVal
        $changed = $changed + 1
    }
}

Write-Host "Updated $changed cells in column G"
